$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the 4th paragraph-of-text ("My role as a mentor...")
# together with the pair of manual line breaks that followed it, so
# the pair of breaks that preceded it now directly separates paragraph
# 3 from the (soon to be replaced) final paragraph.
# ---------------------------------------------------------------------
$old4 = "My role as a mentor and leader in previous positions has involved guiding junior engineers, conducting code reviews, and contributing to design discussions and best practices. I am passionate about staying at the forefront of AI advancements and am committed to leveraging my skills to drive the development and deployment of next-generation AI solutions in your team."

$hit = $d.Content
$hit.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r = $d.Range($hit.Start, $hit.End)
$r.MoveEnd(1, 2)
$r.Text = ""

# ---------------------------------------------------------------------
# Step 2: replace the remaining four paragraphs-of-text with their new
# wording.
# ---------------------------------------------------------------------
$old1 = "I am a senior AI/ML Engineer with over 10 years of experience in the software development industry, specializing in full-stack engineering, cloud solutions, and agentic intelligence. My extensive background in both fintech and healthcare has equipped me with the skills to drive AI initiatives that significantly enhance business operations and decision-making processes. My experience aligns perfectly with the role of Gen AI/Agent AI Engineer, as I have a proven track record of designing, developing, and deploying AI solutions using large language models and agentic frameworks such as LangGraph and AutoGen."
$new1 = "I am a senior AI/ML Engineer with over a decade of experience in the software development industry, specializing in full-stack AI solutions and cloud architecture. My career has been marked by a commitment to leveraging advanced AI and machine learning techniques to drive business impact, particularly in fintech and healthcare domains. I am particularly interested in the opportunity at Mercor as it aligns with my expertise in translating complex datasets into actionable insights using Python, SQL, and data visualization tools."

$old2 = "During my tenure at InsoftAI and CoreWeave, I led projects that involved building and optimizing retrieval-augmented generation pipelines, fine-tuning LLMs with techniques like LoRA and QLoRA, and integrating AI solutions with cloud-native services across AWS, Azure, and GCP. My expertise in working with unstructured data and developing multi-agent systems in the fintech and healthcare sectors has honed my ability to deliver high-impact AI applications that streamline complex workflows and improve operational efficiency."
$new2 = "Throughout my career, I have consistently demonstrated my ability to lead AI initiatives that bridge the gap between technical capabilities and business needs. My experience includes architecting AI-powered credit decisioning systems and healthcare diagnostic platforms, where I applied predictive analytics and multi-agent systems to enhance operational efficiency and accuracy. These projects have honed my skills in statistical modeling, data wrangling, and the development of predictive tools, aligning closely with the responsibilities outlined in the job description."

$old3 = "I possess strong Python programming skills and have hands-on experience with ML/AI libraries such as Hugging Face Transformers and PyTorch. My familiarity with vector databases, including Pinecone and Weaviate, and my proficiency in REST API development and containerization, make me adept at building scalable AI systems that meet the demands of fast-paced, innovation-driven environments."
$new3 = "I am proficient in Python and SQL, with extensive experience in using data visualization tools like Tableau and Power BI to create compelling dashboards that facilitate informed decision-making. My strong communication skills have been instrumental in presenting technical findings and recommendations to non-technical stakeholders, ensuring that strategic insights are effectively conveyed and implemented."

$old5 = "I am eager to bring my expertise in Gen AI and agentic frameworks to your organization, where I can continue to innovate and lead projects that push the boundaries of what AI can achieve. Thank you for considering my application. I look forward to the opportunity to discuss how I can contribute to your team's success."
$new5 = "The opportunity to work with Mercor on a contract basis excites me as it offers the chance to apply my skills in a dynamic and impactful way. I am confident that my background in AI/ML engineering, coupled with my experience in data science and analytics, will enable me to contribute effectively to your strategic analytics initiative. I look forward to the possibility of discussing how I can support Mercor in achieving its goals through advanced data science techniques."

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
